$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999204515344,
    0.99446858115042924,
    0.97217357607617094,
    0.96160852361365012,
    0.95143387749850439,
    0.92692537464003422,
    0.92309135930935371,
    0.91854335783638807,
    0.91447172084687911,
    0.91133109141027679,
    0.91097662550670222,
    0.91055812027810612,
    0.91028081731759158,
    0.91091797263096219,
    0.91289460289243773,
    0.9103881775620124,
    0.90668032968214618,
    0.90557142917054434,
    0.99239231836024389,
    0.97585664595329213,
    0.97445816474062075,
    0.97319365899692456,
    0.96513794162432331,
    0.95211620644786599,
    0.94565907484779066,
    0.92401802942390887,
    0.92021363376161913,
    0.90388847479194201,
    0.89242846943473564,
    0.88763500160232145,
    0.87998127362606693,
    0.87830197734763915,
    0.87778197832358562
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
